$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = '''26.685.27'
$ws.Range("E2").Value = '''  -1.42%  '
$ws.Range("D3").Value = '''1.594.83'
$ws.Range("E3").Value = '''  -1.75%  '
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("D5").Value = '''211.43'
$ws.Range("E5").Value = '''  -1.19%  '
$ws.Range("D6").Value = '''0.511'
$ws.Range("E6").Value = '''  -0.59%  '
$ws.Range("E7").Value = '''  +0.05%  '
$ws.Range("D8").Value = '''0.0617'
$ws.Range("E8").Value = '''  -1.60%  '
$ws.Range("E9").Value = '''  -1.49%  '
$ws.Range("D10").Value = '''19.64'
$ws.Range("E10").Value = '''  -1.17%  '
$ws.Range("D11").Value = '''0.0835'
$ws.Range("E11").Value = '''  -1.21%  '
$ws.Range("D12").Value = '''1.820.09'
$ws.Range("E12").Value = '''  -1.70%  '
$ws.Range("D13").Value = '''1.600.72'
$ws.Range("E13").Value = '''  -1.22%  '
$ws.Range("E14").Value = '''  -2.32%  '
$ws.Range("D15").Value = '''0.523'
$ws.Range("E15").Value = '''  -2.87%  '
$ws.Range("D16").Value = '''64.72'
$ws.Range("E16").Value = '''  +0.57%  '
$ws.Range("D17").Value = '''26.660.53'
$ws.Range("E17").Value = '''  -1.49%  '
$ws.Range("D18").Value = '''0.0₃0730'
$ws.Range("E18").Value = '''  -0.83%  '
$ws.Range("D19").Value = '''208.63'
$ws.Range("E19").Value = '''  -2.69%  '
$ws.Range("E20").Value = '''  +0.04%  '
$ws.Range("D21").Value = '''6.69'
$ws.Range("E21").Value = '''  -1.97%  '
$ws.Range("D22").Value = '''4.24'
$ws.Range("E22").Value = '''  -2.30%  '
$ws.Range("D23").Value = '''2.33'
$ws.Range("E23").Value = '''  +0.07%  '
$ws.Range("D24").Value = '''8.88'
$ws.Range("E24").Value = '''  -1.53%  '
$ws.Range("D25").Value = '''146.62'
$ws.Range("E25").Value = '''  -1.03%  '
$ws.Range("E26").Value = '''  -0.04%  '
$ws.Range("E27").Value = '''  -2.77%  '
$ws.Range("D28").Value = '''0.115'
$ws.Range("E28").Value = '''  +0.09%  '
$ws.Range("D29").Value = '''15.32'
$ws.Range("E29").Value = '''  -1.19%  '
$ws.Range("E30").Value = '''  -1.28%  '
$ws.Range("D31").Value = '''1.15'
$ws.Range("E31").Value = '''  -1.20%  '
$ws.Range("E32").Value = '''  -3.32%  '
$ws.Range("D33").Value = '''0.668'
$ws.Range("E33").Value = '''  -6.20%  '
$ws.Range("D34").Value = '''2.92'
$ws.Range("E34").Value = '''  -2.32%  '
$ws.Range("D35").Value = '''1.294.74'
$ws.Range("E35").Value = '''  -4.35%  '
$ws.Range("E36").Value = '''  -0.62%  '
$ws.Range("E37").Value = '''  -5.31%  '
$ws.Range("D38").Value = '''0.0171'
$ws.Range("E38").Value = '''  -2.95%  '
$ws.Range("D39").Value = '''0.834'
$ws.Range("E39").Value = '''  -0.81%  '
$ws.Range("E40").Value = '''  +0.07%  '
$ws.Range("D41").Value = '''0.795'
$ws.Range("E41").Value = '''  -0.23%  '
$ws.Range("D42").Value = '''2.19'
$ws.Range("E42").Value = '''  -1.66%  '
$ws.Range("D43").Value = '''5.35'
$ws.Range("E43").Value = '''  +0.09%  '
$ws.Range("D44").Value = '''63.49'
$ws.Range("E44").Value = '''  -0.74%  '
$ws.Range("D45").Value = '''1.732.27'
$ws.Range("E45").Value = '''  -1.75%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''0.896'
$ws.Range("E46").Value = '''  +5.43%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''89.80'
$ws.Range("E47").Value = '''  -0.08%  '
$ws.Range("D48").Value = '''1.63'
$ws.Range("E48").Value = '''  -1.35%  '
$ws.Range("D49").Value = '''0.0983'
$ws.Range("E49").Value = '''  -1.95%  '
$ws.Range("E50").Value = '''  -1.72%  '
$ws.Range("D51").Value = '''7.51'
$ws.Range("E51").Value = '''  -1.50%  '
